# Rename the three logo inline pictures embedded in the headers/footers.
#
#   footer (Primary / "default")   -> footer2.xml : image1.png -> image2.png
#   footer (FirstPage / "first")   -> footer1.xml : image1.png -> image2.png
#   header (FirstPage / "first")   -> header1.xml : image2.jpg -> image1.jpg
#
# Word's InlineShape object has no writable .Name property (only the
# floating Shape object does), so each picture is temporarily converted
# to a (floating) Shape, renamed, and converted back to an inline shape -
# this mirrors the <wp:docPr name="..."/> rename that the document diff
# records, while leaving the picture inline exactly as it was.

$d = $word.ActiveDocument

function Rename-InlineLogo($inlineShape, $newName) {
    $shape = $inlineShape.ConvertToShape()
    $shape.Name = $newName
    [void]$shape.ConvertToInlineShape()
}

$sec = $d.Sections.Item(1)

# Footer, Primary ("default") header/footer -> footer2.xml (id="2")
$ftrPrimary = $sec.Footers.Item(1)
Rename-InlineLogo $ftrPrimary.Range.InlineShapes.Item(1) "image2.png"

# Footer, FirstPage ("first") header/footer -> footer1.xml (id="3")
$ftrFirst = $sec.Footers.Item(2)
Rename-InlineLogo $ftrFirst.Range.InlineShapes.Item(1) "image2.png"

# Header, FirstPage ("first") header/footer -> header1.xml (id="1")
$hdrFirst = $sec.Headers.Item(2)
Rename-InlineLogo $hdrFirst.Range.InlineShapes.Item(1) "image1.jpg"

Write-Output "Renamed footer/header logo inline shapes."
